$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price / volume figures scraped on 2024-05-23
$ws.Range("D2").Value = '67.997.94'
$ws.Range("E2").Value = '  -3.28%  '
$ws.Range("D3").Value = '3.810.43'
$ws.Range("E3").Value = '  +1.14%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.49%  '
$ws.Range("D7").Value = '3.809.55'
$ws.Range("E7").Value = '  +1.14%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.62%  '
$ws.Range("E10").Value = '  -4.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.29'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.43%  '
$ws.Range("E12").Value = '  -4.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.03'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.43%  '
$ws.Range("E14").Value = '  -5.02%  '
$ws.Range("D15").Value = '4.443.48'
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("D16").Value = '3.807.09'
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("D17").Value = '68.066.66'
$ws.Range("E17").Value = '  -3.30%  '
$ws.Range("E18").Value = '  -4.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '489.90'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.734'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.76'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -9.36%  '
$ws.Range("E26").Value = '  +2.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.01%  '
$ws.Range("E28").Value = '  -10.19%  '
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.72%  '
$ws.Range("E31").Value = '  -2.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.77'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.96%  '
$ws.Range("E34").Value = '  -4.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.98%  '
$ws.Range("E37").Value = '  -1.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.327'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '449.79'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.20%  '
$ws.Range("E41").Value = '  -2.21%  '
$ws.Range("E42").Value = '  -4.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.90'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.59'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.37%  '
$ws.Range("D46").Value = '2.833.86'
$ws.Range("E46").Value = '  -5.24%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '139.22'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0352'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.24'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.15%  '
$ws.Range("E51").Value = '  -7.11%  '
